# Scheduled market-data refresh: update currentAveragePrice* / Leve*Profit* columns
# across the crafting-class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1891.931
$ws.Range("I70").Value = 1868.3636
$ws.Range("J70").Value = 1966
$ws.Range("K70").Value = 5605.0908
$ws.Range("L70").Value = 5898
$ws.Range("M70").Value = -5335.0908
$ws.Range("N70").Value = -6438

$ws.Range("H73").Value = 1891.931
$ws.Range("I73").Value = 1868.3636
$ws.Range("J73").Value = 1966
$ws.Range("K73").Value = 5605.0908
$ws.Range("L73").Value = 5898
$ws.Range("M73").Value = -4669.0908
$ws.Range("N73").Value = -7770

$ws.Range("H116").Value = 1950.15
$ws.Range("I116").Value = 1901
$ws.Range("J116").Value = 1966.5333
$ws.Range("K116").Value = 1901
$ws.Range("L116").Value = 1966.5333
$ws.Range("M116").Value = 1541
$ws.Range("N116").Value = -8850.533299999999

$ws.Range("H127").Value = 1661.8857
$ws.Range("I127").Value = 568.4286
$ws.Range("J127").Value = 1935.25
$ws.Range("K127").Value = 1705.2858
$ws.Range("L127").Value = 5805.75
$ws.Range("M127").Value = 3254.7142
$ws.Range("N127").Value = -15725.75

$ws.Range("H132").Value = 2519.2942
$ws.Range("I132").Value = 1655.2
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 4965.6
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -2435.6
$ws.Range("N132").Value = -32060

$ws.Range("H137").Value = 1808.2632
$ws.Range("I137").Value = 1104.3846
$ws.Range("J137").Value = 3333.3333
$ws.Range("K137").Value = 3313.1538
$ws.Range("L137").Value = 9999.999899999999
$ws.Range("M137").Value = -763.1538
$ws.Range("N137").Value = -15099.9999

$ws.Range("H138").Value = 4019.7964
$ws.Range("I138").Value = 2053.5
$ws.Range("J138").Value = 4581.595
$ws.Range("K138").Value = 6160.5
$ws.Range("L138").Value = 13744.785
$ws.Range("M138").Value = -1020.5
$ws.Range("N138").Value = -24024.785

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1833.5714
$ws.Range("I2").Value = 1899.2307
$ws.Range("J2").Value = 980
$ws.Range("K2").Value = 1899.2307
$ws.Range("L2").Value = 980
$ws.Range("M2").Value = -1786.2307
$ws.Range("N2").Value = -1206

$ws.Range("H32").Value = 4825.649
$ws.Range("I32").Value = 4511.85
$ws.Range("K32").Value = 4511.85
$ws.Range("M32").Value = -4224.85

$ws.Range("H74").Value = 1681.9656
$ws.Range("I74").Value = 1451.1177
$ws.Range("J74").Value = 2009
$ws.Range("K74").Value = 1451.1177
$ws.Range("L74").Value = 2009
$ws.Range("M74").Value = -577.1177
$ws.Range("N74").Value = -3757

$ws.Range("H77").Value = 1681.9656
$ws.Range("I77").Value = 1451.1177
$ws.Range("J77").Value = 2009
$ws.Range("K77").Value = 7255.5885
$ws.Range("L77").Value = 10045
$ws.Range("M77").Value = -2887.5885
$ws.Range("N77").Value = -18781

$ws.Range("H116").Value = 1833.5714
$ws.Range("I116").Value = 1899.2307
$ws.Range("J116").Value = 980
$ws.Range("K116").Value = 1899.2307
$ws.Range("L116").Value = 980
$ws.Range("M116").Value = 394.7692999999999
$ws.Range("N116").Value = -5568

$ws.Range("H132").Value = 2498.6978
$ws.Range("I132").Value = 2119.6667
$ws.Range("J132").Value = 6194.25
$ws.Range("K132").Value = 6359.000100000001
$ws.Range("L132").Value = 18582.75
$ws.Range("M132").Value = -3829.000100000001
$ws.Range("N132").Value = -23642.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1833.5714
$ws.Range("I3").Value = 1899.2307
$ws.Range("J3").Value = 980
$ws.Range("K3").Value = 1899.2307
$ws.Range("L3").Value = 980
$ws.Range("M3").Value = -1785.2307
$ws.Range("N3").Value = -1208

$ws.Range("H94").Value = 2338.7058
$ws.Range("I94").Value = 1638.625
$ws.Range("J94").Value = 2961
$ws.Range("K94").Value = 1638.625
$ws.Range("L94").Value = 2961
$ws.Range("M94").Value = -1187.625
$ws.Range("N94").Value = -3863

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1290
$ws.Range("I4").Value = 987.5
$ws.Range("K4").Value = 987.5
$ws.Range("M4").Value = -875.5

$ws.Range("H31").Value = 2620.724
$ws.Range("J31").Value = 3825.2856
$ws.Range("L31").Value = 3825.2856
$ws.Range("N31").Value = -4415.2856

$ws.Range("H34").Value = 2620.724
$ws.Range("J34").Value = 3825.2856
$ws.Range("L34").Value = 3825.2856
$ws.Range("N34").Value = -4229.2856

$ws.Range("H58").Value = 1682.9333
$ws.Range("I58").Value = 1107.1111
$ws.Range("K58").Value = 1107.1111
$ws.Range("M58").Value = -904.1111000000001

$ws.Range("H132").Value = 1792.3572
$ws.Range("I132").Value = 1586.25
$ws.Range("J132").Value = 3029
$ws.Range("K132").Value = 4758.75
$ws.Range("L132").Value = 9087
$ws.Range("M132").Value = -2228.75
$ws.Range("N132").Value = -14147

$ws.Range("H136").Value = 1682.9333
$ws.Range("I136").Value = 1107.1111
$ws.Range("K136").Value = 3321.3333
$ws.Range("M136").Value = -771.3333000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7825.3335
$ws.Range("I5").Value = 12756.125
$ws.Range("K5").Value = 38268.375
$ws.Range("M5").Value = -38156.375

$ws.Range("H56").Value = 5166.6665
$ws.Range("I56").Value = 5166.6665
$ws.Range("K56").Value = 5166.6665
$ws.Range("M56").Value = -4636.6665

$ws.Range("H121").Value = 910.1774
$ws.Range("I121").Value = 268
$ws.Range("J121").Value = 1005.3148
$ws.Range("K121").Value = 804
$ws.Range("L121").Value = 3015.9444
$ws.Range("M121").Value = 506
$ws.Range("N121").Value = -5635.9444

$ws.Range("H129").Value = 1496.4706
$ws.Range("I129").Value = 1000.8333
$ws.Range("J129").Value = 1766.8182
$ws.Range("K129").Value = 3002.4999
$ws.Range("L129").Value = 5300.4546
$ws.Range("M129").Value = 1997.5001
$ws.Range("N129").Value = -15300.4546

$ws.Range("H131").Value = 16667726
$ws.Range("I131").Value = 6667342
$ws.Range("J131").Value = 19609016
$ws.Range("K131").Value = 20002026
$ws.Range("L131").Value = 58827048
$ws.Range("M131").Value = -19996986
$ws.Range("N131").Value = -58837128

$ws.Range("H135").Value = 7825.3335
$ws.Range("I135").Value = 12756.125
$ws.Range("K135").Value = 114805.125
$ws.Range("M135").Value = -112270.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 122.125
$ws.Range("I2").Value = 151.75
$ws.Range("K2").Value = 151.75
$ws.Range("M2").Value = -38.75

$ws.Range("H122").Value = 81910504
$ws.Range("I122").Value = 88736136
$ws.Range("J122").Value = 3008
$ws.Range("K122").Value = 266208408
$ws.Range("L122").Value = 9024
$ws.Range("M122").Value = -266205958
$ws.Range("N122").Value = -13924

$ws.Range("H132").Value = 3607.8628
$ws.Range("I132").Value = 3500.1304
$ws.Range("J132").Value = 3696.3572
$ws.Range("K132").Value = 10500.3912
$ws.Range("L132").Value = 11089.0716
$ws.Range("M132").Value = -7970.3912
$ws.Range("N132").Value = -16149.0716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2833.3333
$ws.Range("I7").Value = 2500
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -2388
$ws.Range("N7").Value = -3724

$ws.Range("H122").Value = 8145104.5
$ws.Range("I122").Value = 7939005
$ws.Range("J122").Value = 10000000
$ws.Range("K122").Value = 23817015
$ws.Range("L122").Value = 30000000
$ws.Range("M122").Value = -23814565
$ws.Range("N122").Value = -30004900

$ws.Range("H126").Value = 2833.3333
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -15440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H108").Value = 40416.668
$ws.Range("J108").Value = 40416.668
$ws.Range("L108").Value = 40416.668
$ws.Range("N108").Value = -48096.668

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H126").Value = 834.53845
$ws.Range("I126").Value = 723.1111
$ws.Range("J126").Value = 1085.25
$ws.Range("K126").Value = 2169.3333
$ws.Range("L126").Value = 3255.75
$ws.Range("M126").Value = 300.6667000000002
$ws.Range("N126").Value = -8195.75
